# #355 nested struct implemented
# Adds two new worksheets ("NestedStructure" and "StructureOrder") describing
# ExcelSchemaBuilder test data for nested struct fields, after the existing
# "Stringfields-All" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet: NestedStructure
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsNested = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsNested.Name = "NestedStructure"

$wsNested.Cells.Item(1, 1).Value = "class"
$wsNested.Cells.Item(1, 2).Value = "name"
$wsNested.Cells.Item(1, 3).Value = "useSequence"
$wsNested.Cells.Item(1, 4).Value = "multiplicity"
$wsNested.Range("A1:D1").Font.Bold = $true

$nestedRows = @(
    @("struct", "NestedStructure", $true, $null),
    @("field",  "stringField1",    $null, $null),
    @("struct", "SubStruct",       $true, "0..*"),
    @("field",  "stringField11",   $null, $null),
    @("struct", "SubStruct",       $null, $null),
    @("struct", "NestedStructure", $null, $null)
)

$r = 2
foreach ($row in $nestedRows) {
    $wsNested.Cells.Item($r, 1).Value = $row[0]
    $wsNested.Cells.Item($r, 2).Value = $row[1]
    if ($null -ne $row[2]) { $wsNested.Cells.Item($r, 3).Value = $row[2] }
    if ($null -ne $row[3]) { $wsNested.Cells.Item($r, 4).Value = $row[3] }
    $r++
}

$wsNested.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait

$wsNested.Activate() | Out-Null
$wsNested.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$wsNested.Range("C1").Select() | Out-Null

# ---------------------------------------------------------------------------
# Sheet: StructureOrder
# ---------------------------------------------------------------------------
$wsOrder = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsNested)
$wsOrder.Name = "StructureOrder"

$wsOrder.Cells.Item(1, 1).Value = "class"
$wsOrder.Cells.Item(1, 2).Value = "name"
$wsOrder.Range("A1:B1").Font.Bold = $true

$orderRows = @(
    @("struct", "StructureOrder"),
    @("field",  "stringField1"),
    @("struct", "SubStruct1"),
    @("field",  "stringField11"),
    @("struct", "SubStruct1"),
    @("field",  "stringField2"),
    @("struct", "SubStruct2"),
    @("field",  "stringField21"),
    @("struct", "SubStruct2"),
    @("field",  "stringField3"),
    @("struct", "StructureOrder")
)

$r = 2
foreach ($row in $orderRows) {
    $wsOrder.Cells.Item($r, 1).Value = $row[0]
    $wsOrder.Cells.Item($r, 2).Value = $row[1]
    $r++
}

$wsOrder.PageSetup.Orientation = [Microsoft.Office.Interop.Excel.XlPageOrientation]::xlPortrait

$wsOrder.Activate() | Out-Null
$wsOrder.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$wsOrder.Range("B11").Select() | Out-Null
